# Doing Updates for Financials
# Apply value updates to the CPA yearly financials sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Income Statement ---------------------------------------------------

# Total Revenue (row 8)
$ws.Range("D8").Value = 2521800

# Cost of Revenue (row 9)
$ws.Range("D9").Value = 1182200
$ws.Range("F9").Value = 1126600
$ws.Range("G9").Value = 1187100
$ws.Range("H9").Value = 1120300
$ws.Range("I9").Value = 1026000
$ws.Range("J9").Value = 784300

# Gross Profit (row 10)
$ws.Range("D10").Value = 1339600
$ws.Range("F10").Value = 1127100
$ws.Range("G10").Value = 1518000
$ws.Range("H10").Value = 1488000
$ws.Range("I10").Value = 1223400
$ws.Range("J10").Value = 1046600

# Non Recurring (row 14): D14 was "NA" text, now numeric 0
$ws.Range("D14").Value = 0

# Others (row 15)
$ws.Range("D15").Value = 167300
$ws.Range("E15").Value = "NA"
$ws.Range("F15").Value = 134900
$ws.Range("G15").Value = 115100
$ws.Range("H15").Value = 137400
$ws.Range("I15").Value = 89200
$ws.Range("J15").Value = 75500

# Total Operating Expenses (row 17)
$ws.Range("D17").Value = 2097800

# Operating Income or Loss (row 18)
$ws.Range("D18").Value = 424000

# Total Other Income/Expenses Net (row 20)
$ws.Range("D20").Value = 24500

# Earnings Before Interest And Taxes (row 21)
$ws.Range("D21").Value = 612900
$ws.Range("J21").Value = "NA"

# Income Before Tax (row 23)
$ws.Range("D23").Value = 413400

# Income Tax Expense (row 24)
$ws.Range("D24").Value = 49300

# Income After Tax (row 26)
$ws.Range("D26").Value = 364000

# Net Income From Continuing Ops (row 27)
$ws.Range("D27").Value = 364000

# Other Items (row 32)
$ws.Range("D32").Value = -24500

# Net Income (row 33)
$ws.Range("D33").Value = 364000

# Net Income Applicable To Common Shares (row 35)
$ws.Range("D35").Value = 364000

# --- Balance Sheet --------------------------------------------------------

# Property Plant and Equipment (row 48)
$ws.Range("D48").Value = 5443300

# Other Assets (row 52)
$ws.Range("D52").Value = 96100

# Total Assets (row 54)
$ws.Range("D54").Value = 4045000

# Accounts Payable (row 57)
$ws.Range("D57").Value = 212400

# Other Current Liabilities (row 59)
$ws.Range("D59").Value = 621000

# Total Current Liabilities (row 60)
$ws.Range("D60").Value = 1057500

# Other Liabilities (row 62)
$ws.Range("D62").Value = 216200

# Total Liabilities (row 66)
$ws.Range("D66").Value = 2149800

# Retained Earnings (row 72)
$ws.Range("D72").Value = 1934000

# Total Stockholder Equity (row 76)
$ws.Range("D76").Value = 1895100

# --- Cash Flow Statement ---------------------------------------------------

# Net Income (row 81)
$ws.Range("D81").Value = 364000

# Depreciation (row 83)
$ws.Range("J83").Value = "NA"

# Capital Expenditures (row 91)
$ws.Range("D91").Value = -109900
$ws.Range("E91").Value = -88300
$ws.Range("F91").Value = -81800
$ws.Range("G91").Value = -106300

# Other Cashflows from Investing Activities (row 94)
$ws.Range("J94").Value = "NA"

# Net Borrowings (row 100)
$ws.Range("J100").Value = "NA"

# Other Cash Flows from Financing Activities (row 101)
$ws.Range("J101").Value = "NA"
